$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("hus"): update monthly values from 4560 to 1000
$ws.Range("B4:M4").Value = 1000

# Row 12 ("transport"): update values from 23 to 120 (only B, D, F, H, J, L have values)
$ws.Range("B12").Value = 120
$ws.Range("D12").Value = 120
$ws.Range("F12").Value = 120
$ws.Range("H12").Value = 120
$ws.Range("J12").Value = 120
$ws.Range("L12").Value = 120

# Row 13 ("overført fra"): update values from 1500 to 1200 (only C, F, I, L have values)
$ws.Range("C13").Value = 1200
$ws.Range("F13").Value = 1200
$ws.Range("I13").Value = 1200
$ws.Range("L13").Value = 1200

# Force recalculation so dependent formula cells (row 5, N4, N12, N13, row 14) update
$excel.Calculate()

# Update the selected cell on the sheet view
$ws.Range("K14").Select()
